# update scripts wuth new tpm
# Recomputed NATMI ligand-receptor metrics (columns G:T) for Wnt4-Fzd6 edges
# using the new TPM-based expression values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.4173876666666667
$ws.Range("H2").Value = 1.252163
$ws.Range("I2").Value = 0.1865415014963835
$ws.Range("J2").Value = 0.1865415014963835
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.502639
$ws.Range("N2").Value = 31.507917
$ws.Range("O2").Value = 0.9701256668284471
$ws.Range("P2").Value = 0.970125666828447
$ws.Range("Q2").Value = 4.383671986052334
$ws.Range("R2").Value = 39.453047874471
$ws.Range("S2").Value = 0.1809686985303589
$ws.Range("T2").Value = 0.1809686985303588

# Row 3
$ws.Range("G3").Value = 0.4173876666666667
$ws.Range("H3").Value = 1.252163
$ws.Range("I3").Value = 0.1865415014963835
$ws.Range("J3").Value = 0.1865415014963835
$ws.Range("O3").Value = 0.02703852164627077
$ws.Range("P3").Value = 0.02703852164627077
$ws.Range("Q3").Value = 0.1221779960451111
$ws.Range("R3").Value = 1.099601964406
$ws.Range("S3").Value = 0.005043806426137818
$ws.Range("T3").Value = 0.005043806426137817

# Row 4
$ws.Range("G4").Value = 0.4173876666666667
$ws.Range("H4").Value = 1.252163
$ws.Range("I4").Value = 0.1865415014963835
$ws.Range("J4").Value = 0.1865415014963835
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03070066666666667
$ws.Range("N4").Value = 0.092102
$ws.Range("O4").Value = 0.002835811525282158
$ws.Range("P4").Value = 0.002835811525282158
$ws.Range("Q4").Value = 0.01281407962511111
$ws.Range("R4").Value = 0.115326716626
$ws.Range("S4").Value = 0.0005289965398868835
$ws.Range("T4").Value = 0.0005289965398868835

# Row 5
$ws.Range("I5").Value = 0.5456214114616024
$ws.Range("J5").Value = 0.5456214114616023
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 10.502639
$ws.Range("N5").Value = 31.507917
$ws.Range("O5").Value = 0.9701256668284471
$ws.Range("P5").Value = 0.970125666828447
$ws.Range("Q5").Value = 12.821947273009
$ws.Range("R5").Value = 115.397525457081
$ws.Range("S5").Value = 0.5293213356300656
$ws.Range("T5").Value = 0.5293213356300653

# Row 6
$ws.Range("I6").Value = 0.5456214114616024
$ws.Range("J6").Value = 0.5456214114616023
$ws.Range("O6").Value = 0.02703852164627077
$ws.Range("P6").Value = 0.02703852164627077
$ws.Range("R6").Value = 3.216262177866001
$ws.Range("S6").Value = 0.01475279634447335
$ws.Range("T6").Value = 0.01475279634447334

# Row 7
$ws.Range("I7").Value = 0.5456214114616024
$ws.Range("J7").Value = 0.5456214114616023
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.03070066666666667
$ws.Range("N7").Value = 0.092102
$ws.Range("O7").Value = 0.002835811525282158
$ws.Range("P7").Value = 0.002835811525282158
$ws.Range("Q7").Value = 0.03748032558733334
$ws.Range("R7").Value = 0.3373229302860001
$ws.Range("S7").Value = 0.001547279487063531
$ws.Range("T7").Value = 0.00154727948706353

# Row 8
$ws.Range("G8").Value = 0.599287
$ws.Range("H8").Value = 1.797861
$ws.Range("I8").Value = 0.2678370870420142
$ws.Range("J8").Value = 0.2678370870420142
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 10.502639
$ws.Range("N8").Value = 31.507917
$ws.Range("O8").Value = 0.9701256668284471
$ws.Range("P8").Value = 0.970125666828447
$ws.Range("Q8").Value = 6.294095018393
$ws.Range("R8").Value = 56.646855165537
$ws.Range("S8").Value = 0.2598356326680228
$ws.Range("T8").Value = 0.2598356326680228

# Row 9
$ws.Range("G9").Value = 0.599287
$ws.Range("H9").Value = 1.797861
$ws.Range("I9").Value = 0.2678370870420142
$ws.Range("J9").Value = 0.2678370870420142
$ws.Range("O9").Value = 0.02703852164627077
$ws.Range("P9").Value = 0.02703852164627077
$ws.Range("Q9").Value = 0.1754236901646667
$ws.Range("R9").Value = 1.578813211482
$ws.Range("S9").Value = 0.00724191887565961
$ws.Range("T9").Value = 0.007241918875659609

# Row 10
$ws.Range("G10").Value = 0.599287
$ws.Range("H10").Value = 1.797861
$ws.Range("I10").Value = 0.2678370870420142
$ws.Range("J10").Value = 0.2678370870420142
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.03070066666666667
$ws.Range("N10").Value = 0.092102
$ws.Range("O10").Value = 0.002835811525282158
$ws.Range("P10").Value = 0.002835811525282158
$ws.Range("Q10").Value = 0.01839851042466667
$ws.Range("R10").Value = 0.165586593822
$ws.Range("S10").Value = 0.0007595354983317445
$ws.Range("T10").Value = 0.0007595354983317445
